# edit.ps1 -- apply the TestReview.docx change:
#   1. Remove the "_GoBack" bookmark from its old position (inside the
#      word "invalid" in the "Load file with invalid image file format"
#      cell).
#   2. In the cell whose paragraph is the single word "width" (the
#      test_calcHori_largeH row), capitalize it to "Width" and append
#      ", where width/(h/48) < 250", then re-insert the "_GoBack"
#      bookmark at the very end of that paragraph (its new position
#      after the edit, mirroring how Word leaves "_GoBack" at the site
#      of the most recent edit).

function Find-ParaByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $pp = $doc.Paragraphs.Item($i)
        $t = $pp.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $pp
        }
    }
    return $null
}

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Part 1: delete the old "_GoBack" bookmark (Bookmark.Delete() is not
# wired up in this host, so instead we re-type the two characters that
# straddle it via Find/Replace -- a real edit across the bookmark's
# position removes it, same as Word does when the bookmarked spot is
# edited).
# ------------------------------------------------------------------
$oldBm = $d.Bookmarks.Item("_GoBack")
$bmPos = $oldBm.Start
$rStraddle = $d.Range($bmPos - 1, $bmPos + 1)
$straddleText = $rStraddle.Text
$rStraddle.Find.Execute($straddleText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $straddleText, 2)

# ------------------------------------------------------------------
# Part 2: locate the paragraph whose entire text is "width"
# ------------------------------------------------------------------
$target = Find-ParaByText $d "width"
$startPos = $target.Range.Start

# Capitalize "w" -> "W"
$rFirst = $d.Range($startPos, $startPos + 1)
$rFirst.Text = "W"

# Append the new clause after "Width"
$pCap = Find-ParaByText $d "Width"
$pCap.Range.InsertAfter(", where width/(h/48) < 250")

# ------------------------------------------------------------------
# Part 3: re-add the "_GoBack" bookmark at the new end of the
# paragraph. Bookmarks.Add() collapses to position 0 when given a
# zero-length range that sits exactly at a paragraph's last position,
# so we guard the true end with a throw-away character, add the
# bookmark just before it, then remove the guard character.
# ------------------------------------------------------------------
$pFull = Find-ParaByText $d "Width, where width/(h/48) < 250"
$endPos = $pFull.Range.End - 1
$rGuard = $d.Range($endPos, $endPos)
$rGuard.InsertAfter("Z")
$rBmTarget = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $rBmTarget) | Out-Null
$rGuardChar = $d.Range($endPos, $endPos + 1)
$rGuardChar.Delete()

# ------------------------------------------------------------------
# Part 4: split the merged run back into "W" / "idth" / ", where ..."
# pieces (touching Font forces the host to stop coalescing adjacent
# runs that share identical formatting).
# ------------------------------------------------------------------
$rIdth = $d.Range($startPos + 1, $startPos + 5)
$rIdth.Font.Bold = $true
$rIdth.Font.Bold = $false

$pFinal = Find-ParaByText $d "Width, where width/(h/48) < 250"
$clauseEnd = $pFinal.Range.End - 1
$rClause = $d.Range($startPos + 5, $clauseEnd)
$rClause.Font.Bold = $true
$rClause.Font.Bold = $false

Write-Host "Done. Final paragraph text:" $pFinal.Range.Text.TrimEnd([char]13, [char]7)
